$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.675.51"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.059.41"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "516.45"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "'139.80"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.433"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").Value = "7.27"
$ws.Range("E9").Value = "  -3.82%  "
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "3.576.51"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("E13").Value = "  -3.31%  "
$ws.Range("D14").Value = "26.75"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "57.622.56"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "6.21"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "3.069.26"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").Value = "13.33"
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").Value = "8.16"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").Value = "328.86"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "0.507"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "66.01"
$ws.Range("E24").Value = "  +2.19%  "
$ws.Range("D25").Value = "3.175.72"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "0.0₃0902"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").Value = "6.67"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("D33").Value = "20.79"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").Value = "154.04"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "4.61"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").Value = "25.42"
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").Value = "0.0676"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "37.11"
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").Value = "3.87"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.669"
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "2.204.11"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").Value = "'1.40"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").Value = "'6.10"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("D49").Value = "20.04"
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "1.77"
$ws.Range("E50").Value = "  -5.49%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "0.184"
$ws.Range("E51").Value = "  -0.34%  "
